$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F21").Value = "wear:Blades"
$ws.Range("F22").Value = "wear:Blades"
$ws.Range("F23").Value = "wear:Blades"
$ws.Range("F24").Value = "wear:Tower, wear:Drivetrain"
$ws.Range("F25").Value = "wear:Tower, wear:Drivetrain"
$ws.Range("F26").Value = "wear:Tower"
$ws.Range("F31").Value = "wear:TransformerAndInductor"
$ws.Range("F32").Value = "wear:TransformerAndInductor"
$ws.Range("F33").Value = "wear:TransformerAndInductor"
$ws.Range("F34").Value = "wear:TransformerAndInductor"
$ws.Range("F35").Value = "wear:TransformerAndInductor"
$ws.Range("F36").Value = "wear:TransformerAndInductor"
$ws.Range("F37").Value = "wear:TransformerAndInductor"
$ws.Range("F38").Value = "wear:TransformerAndInductor"
$ws.Range("F39").Value = "wear:TransformerAndInductor"
$ws.Range("F40").Value = "wear:TransformerAndInductor"
$ws.Range("F41").Value = "wear:TransformerAndInductor"
